$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new metric data row (row 41) with the latest timestamp reading.
$ws.Cells.Item(41, 1).Value = "2025-04-29 03:57:37"
$ws.Cells.Item(41, 2).Value = 103
